# Commit: "added some more links I just found"
#
# A new slide ("Even More Links") is inserted right after slide 17
# ("For More Information") and before the closing "Thank You!" slide,
# which pushes "Thank You!" down to become the new last slide.

$p = $ppt.ActivePresentation

# "Title and Content" is custom layout #2 on the slide master - the same
# layout already used by the neighboring slides.
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert the new slide at position 18 (pushes the old slide 18 "Thank You!"
# down to position 19).
$slide = $p.Slides.AddSlide(18, $layout)

# --- Title -------------------------------------------------------------
$slide.Shapes.Item(1).TextFrame.TextRange.Text = "Even More Links"

# --- Body / content placeholder -----------------------------------------
$body = $slide.Shapes.Item(2).TextFrame.TextRange

$body.Text = "VSLive: Hate JavaScript? Try TypeScript!" + "`r" + `
             "By Ben Hoelting" + "`r" + `
             "http://visualstudiomagazine.com/blogs/vs-live-video/2013/05/typescript-video.aspx" + "`r" + `
             "A TypeScript Primer" + "`r" + `
             "By Mark Michaelis" + "`r" + `
             "http://visualstudiomagazine.com/articles/2013/02/01/typescript.aspx" + "`r" + `
             "Test-Driven Development with TypeScript" + "`r" + `
             "By Peter Vogel" + "`r" + `
             "http://visualstudiomagazine.com/articles/2013/06/01/test-driven-development-with-typescript.aspx"

# Indent level 2 (COM 1-based) == lvl="1" in the OOXML.
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(3, 1).IndentLevel = 2
$body.Paragraphs(5, 1).IndentLevel = 2
$body.Paragraphs(6, 1).IndentLevel = 2
$body.Paragraphs(8, 1).IndentLevel = 2
$body.Paragraphs(9, 1).IndentLevel = 2

# Hyperlink the "VS Live" blog post URL (paragraph 3) - split into the
# "http://" run and the remainder, both pointing at the same address, to
# mirror how the original author's edit session split the run.
$url1 = "http://visualstudiomagazine.com/blogs/vs-live-video/2013/05/typescript-video.aspx"
$p3 = $body.Paragraphs(3, 1)
$p3.Characters(1, 7).ActionSettings(1).Hyperlink.Address = $url1
$p3.Characters(8, $p3.Length - 7).ActionSettings(1).Hyperlink.Address = $url1

# Hyperlink the "TypeScript Primer" article URL (paragraph 6) - same
# run-split treatment.
$url2 = "http://visualstudiomagazine.com/articles/2013/02/01/typescript.aspx"
$p6 = $body.Paragraphs(6, 1)
$p6.Characters(1, 7).ActionSettings(1).Hyperlink.Address = $url2
$p6.Characters(8, $p6.Length - 7).ActionSettings(1).Hyperlink.Address = $url2

# NOTE: the final URL (paragraph 9, the TDD article) is left as plain,
# un-hyperlinked text, matching the source deck.
